$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 691.8
$ws.Range("I12").Value = 614.75
$ws.Range("K12").Value = 614.75
$ws.Range("M12").Value = -444.75
$ws.Range("H19").Value = 151.16667
$ws.Range("I19").Value = 183.75
$ws.Range("J19").Value = 86
$ws.Range("K19").Value = 183.75
$ws.Range("L19").Value = 86
$ws.Range("M19").Value = -8.75
$ws.Range("N19").Value = -436
$ws.Range("H58").Value = 242.5
$ws.Range("I58").Value = 242.5
$ws.Range("K58").Value = 727.5
$ws.Range("M58").Value = -577.5
$ws.Range("H107").Value = 1145.8966
$ws.Range("I107").Value = 845.36
$ws.Range("K107").Value = 845.36
$ws.Range("M107").Value = 1074.64
$ws.Range("H112").Value = 2282.1765
$ws.Range("I112").Value = 1624.5
$ws.Range("J112").Value = 2484.5386
$ws.Range("K112").Value = 4873.5
$ws.Range("L112").Value = 7453.6158
$ws.Range("M112").Value = -3765.5
$ws.Range("N112").Value = -9669.6158
$ws.Range("H132").Value = 4582.793
$ws.Range("I132").Value = 4197.4
$ws.Range("K132").Value = 12592.2
$ws.Range("M132").Value = -10062.2
$ws.Range("H137").Value = 2374.9
$ws.Range("I137").Value = 1968.625
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 5905.875
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -3355.875
$ws.Range("N137").Value = -17100
$ws.Range("H138").Value = 4355.4443
$ws.Range("J138").Value = 4407.615
$ws.Range("L138").Value = 13222.845
$ws.Range("N138").Value = -23502.845

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2299.8125
$ws.Range("I2").Value = 1528.4
$ws.Range("K2").Value = 1528.4
$ws.Range("M2").Value = -1415.4
$ws.Range("H32").Value = 21173.742
$ws.Range("I32").Value = 17116.793
$ws.Range("K32").Value = 17116.793
$ws.Range("M32").Value = -16829.793
$ws.Range("H64").Value = 33333.332
$ws.Range("J64").Value = 33333.332
$ws.Range("L64").Value = 33333.332
$ws.Range("N64").Value = -33829.332
$ws.Range("H67").Value = 33333.332
$ws.Range("J67").Value = 33333.332
$ws.Range("L67").Value = 33333.332
$ws.Range("N67").Value = -35049.332
$ws.Range("H74").Value = 2753.3333
$ws.Range("I74").Value = 2753.3333
$ws.Range("K74").Value = 2753.3333
$ws.Range("M74").Value = -1879.3333
$ws.Range("H77").Value = 2753.3333
$ws.Range("I77").Value = 2753.3333
$ws.Range("K77").Value = 13766.6665
$ws.Range("M77").Value = -9398.666499999999
$ws.Range("H116").Value = 2299.8125
$ws.Range("I116").Value = 1528.4
$ws.Range("K116").Value = 1528.4
$ws.Range("M116").Value = 765.5999999999999
$ws.Range("H122").Value = 799
$ws.Range("I122").Value = 799
$ws.Range("K122").Value = 2397
$ws.Range("M122").Value = 53

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2299.8125
$ws.Range("I3").Value = 1528.4
$ws.Range("K3").Value = 1528.4
$ws.Range("M3").Value = -1414.4
$ws.Range("H20").Value = 6875
$ws.Range("I20").Value = 7666.6665
$ws.Range("J20").Value = 4500
$ws.Range("K20").Value = 7666.6665
$ws.Range("L20").Value = 4500
$ws.Range("M20").Value = -7419.6665
$ws.Range("N20").Value = -4994
$ws.Range("H82").Value = 19527.375
$ws.Range("H85").Value = 19527.375
$ws.Range("H134").Value = 7750
$ws.Range("I134").Value = 15000
$ws.Range("J134").Value = 5333.3335
$ws.Range("K134").Value = 45000
$ws.Range("L134").Value = 16000.0005
$ws.Range("M134").Value = -42465
$ws.Range("N134").Value = -21070.0005
$ws.Range("H138").Value = 60000
$ws.Range("J138").Value = 60000
$ws.Range("L138").Value = 60000
$ws.Range("N138").Value = -70280

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5287.7144
$ws.Range("I31").Value = 4899.6665
$ws.Range("K31").Value = 4899.6665
$ws.Range("M31").Value = -4604.6665
$ws.Range("H34").Value = 5287.7144
$ws.Range("I34").Value = 4899.6665
$ws.Range("K34").Value = 4899.6665
$ws.Range("M34").Value = -4697.6665
$ws.Range("H41").Value = 18700
$ws.Range("J41").Value = 22500
$ws.Range("L41").Value = 22500
$ws.Range("N41").Value = -23356
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H59").Value = 29976.062
$ws.Range("J59").Value = 34995.445
$ws.Range("L59").Value = 34995.445
$ws.Range("N59").Value = -37285.445
$ws.Range("H60").Value = 10937.6
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H132").Value = 1147.2727
$ws.Range("I132").Value = 1147.2727
$ws.Range("K132").Value = 3441.8181
$ws.Range("M132").Value = -911.8181
$ws.Range("H134").Value = 9096.25
$ws.Range("I134").Value = 8253.286
$ws.Range("K134").Value = 24759.858
$ws.Range("M134").Value = -22224.858

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 465.83334
$ws.Range("H17").Value = 305
$ws.Range("I17").Value = 87.5
$ws.Range("J17").Value = 450
$ws.Range("K17").Value = 262.5
$ws.Range("L17").Value = 1350
$ws.Range("M17").Value = -93.5
$ws.Range("N17").Value = -1688
$ws.Range("H107").Value = 1383.375
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1383.375
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 4150.125
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -7990.125
$ws.Range("H122").Value = 2856.5
$ws.Range("I122").Value = 1807
$ws.Range("K122").Value = 16263
$ws.Range("M122").Value = -13813

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1979.1482
$ws.Range("I102").Value = 1856.2273
$ws.Range("J102").Value = 2520
$ws.Range("K102").Value = 1856.2273
$ws.Range("L102").Value = 2520
$ws.Range("M102").Value = -234.2273
$ws.Range("N102").Value = -5764
$ws.Range("H122").Value = 41302.707
$ws.Range("I122").Value = 45126.383
$ws.Range("K122").Value = 135379.149
$ws.Range("M122").Value = -132929.149
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
$ws.Range("H133").Value = 88081.5
$ws.Range("J133").Value = 88081.5
$ws.Range("L133").Value = 88081.5
$ws.Range("N133").Value = -98201.5
$ws.Range("H135").Value = 20000
$ws.Range("I135").Value = 20000
$ws.Range("K135").Value = 20000
$ws.Range("M135").Value = -14930
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3571.4614
$ws.Range("I16").Value = 3660.75
$ws.Range("K16").Value = 3660.75
$ws.Range("M16").Value = -3490.75
$ws.Range("H40").Value = 4140.077
$ws.Range("I40").Value = 4239.25
$ws.Range("J40").Value = 2950
$ws.Range("K40").Value = 4239.25
$ws.Range("L40").Value = 2950
$ws.Range("M40").Value = -4103.25
$ws.Range("N40").Value = -3222
$ws.Range("H132").Value = 18828.938
$ws.Range("I132").Value = 20842.092
$ws.Range("K132").Value = 62526.276
$ws.Range("M132").Value = -59996.276

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3400.1667
$ws.Range("I126").Value = 3080.3
$ws.Range("J126").Value = 4999.5
$ws.Range("K126").Value = 9240.900000000001
$ws.Range("L126").Value = 14998.5
$ws.Range("M126").Value = -6770.900000000001
$ws.Range("N126").Value = -19938.5
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
